$wb = $excel.ActiveWorkbook

# --- Poland: select the existing "MX-DPBX/MX-BBX" reference rows (no data change) ---
$wsPoland = $wb.Worksheets.Item("Poland")
$wsPoland.Range("A9:A10").Select() | Out-Null

# --- Hungary: insert two new rows for the Panel Accessories (MX-DPBX / MX-BBX) ---
$wsHungary = $wb.Worksheets.Item("Hungary")
$wsHungary.Rows.Item(11).Insert() | Out-Null
$wsHungary.Rows.Item(11).Insert() | Out-Null
$wsHungary.Range("A10").Copy() | Out-Null
$wsHungary.Range("A11:A12").PasteSpecial(-4122) | Out-Null
$wsHungary.Range("A11").Value = "MX-DPBX"
$wsHungary.Range("A12").Value = "MX-BBX"
$wsHungary.Range("A11:A12").Select() | Out-Null

# --- Turkey: insert two new rows for the Panel Accessories (MX-DPBX / MX-BBX) ---
$wsTurkey = $wb.Worksheets.Item("Turkey")
$wsTurkey.Rows.Item(11).Insert() | Out-Null
$wsTurkey.Rows.Item(11).Insert() | Out-Null
$wsTurkey.Range("A10").Copy() | Out-Null
$wsTurkey.Range("A11:A12").PasteSpecial(-4122) | Out-Null
$wsTurkey.Range("A11").Value = "MX-DPBX"
$wsTurkey.Range("A12").Value = "MX-BBX"
$wsTurkey.Range("A11:A12").Select() | Out-Null

# --- Spain: insert one new row and overwrite the following two rows ---
$wsSpain = $wb.Worksheets.Item("Spain")
$wsSpain.Rows.Item(11).Insert() | Out-Null
$wsSpain.Range("A10").Copy() | Out-Null
$wsSpain.Range("A11:A12").PasteSpecial(-4122) | Out-Null
$wsSpain.Range("A11").Value = "MX-DPBX"
$wsSpain.Range("A12").Value = "MX-BBX"
$wsSpain.Range("A11:A12").Select() | Out-Null
